# Questions_Journey.xlsx - add "Missing Number In Arithmetic Progression" row
# and repurpose column E from a "real" check-mark column into a "type" (category) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row: column E header changes from "real" to "type"
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "type"

# ---------------------------------------------------------------------------
# 2) The TestDome rows (4-8) no longer carry a value in the (now repurposed)
#    "type" column, so clear column E for those rows.
# ---------------------------------------------------------------------------
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()

# ---------------------------------------------------------------------------
# 3) The existing LeetCode rows (9-11) are tagged with their algorithm type "DP"
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = "DP"
$ws.Range("E10").Value = "DP"
$ws.Range("E11").Value = "DP"

# ---------------------------------------------------------------------------
# 4) Add the new question row (row 12) -
#    MissingNumberInArithmeticProgression / LeetCode / Math / done
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "MissingNumberInArithmeticProgression"
$ws.Range("D12").Value = "LeetCode"
$ws.Range("E12").Value = "Math"
$ws.Range("F12").Value = "https://leetcode-cn.com/problems/missing-number-in-arithmetic-progression/"
$ws.Range("G12").Value = "√"

# Turn F12 into a real hyperlink, like F10/F11
$ws.Hyperlinks.Add($ws.Range("F12"), "https://leetcode-cn.com/problems/missing-number-in-arithmetic-progression/") | Out-Null

# Re-apply the same visual style used by the existing hyperlink cells (F11)
# since adding the hyperlink above creates a slightly different style entry.
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Column widths: column C (title) and column E (type) grow to fit the new,
#    longer content.
# ---------------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 33.357142857142904
$ws.Columns("E:E").ColumnWidth = 4.642857142857082

# ---------------------------------------------------------------------------
# 6) Update the saved cursor/selection position to G13, matching where the
#    user ended up after entering the new row.
# ---------------------------------------------------------------------------
$ws.Range("G13").Select() | Out-Null

Write-Host "Edit complete"
